$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prix Spot": append a new day column (BE) with header "09-aug" and
# the 24 hourly prices, mirroring the formatting already used for column BD.
# ---------------------------------------------------------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

# Copy the header cell's formatting (bold, centered, bordered) onto the new
# header cell before writing its value.
$wsSpot.Range("BD1").Copy()
$wsSpot.Range("BE1").PasteSpecial(-4122)

$newDayHeader = "09-aug"
$hourlyPrices = @(
    86.26000000000001,
    70.13,
    54.77,
    42,
    41.4,
    44.03,
    50.95,
    47.81,
    46.33,
    15.34,
    -0.02,
    -1,
    -0.02,
    -4.74,
    -7.75,
    -3.96,
    -0.05,
    3.52,
    28.61,
    69.98999999999999,
    85.70999999999999,
    78.95,
    86.33,
    78.91
)

$wsSpot.Range("BE1").Value = $newDayHeader

for ($i = 0; $i -lt $hourlyPrices.Length; $i++) {
    $row = $i + 2
    $wsSpot.Cells.Item($row, 57).Value = $hourlyPrices[$i]
}

# ---------------------------------------------------------------------------
# Sheet "Gaz": append a new daily row (54) for 2025-08-07.
# The date column stores plain text ("2025-08-06", ...), not real dates, so
# force text formatting before writing to stop Excel's autodetect from
# turning the string into a date serial number; then restore the default
# "Normal" style so the cell matches its unstyled neighbours.
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
$gazDateCell = $wsGaz.Range("A54")
$gazDateCell.NumberFormat = "@"
$gazDateCell.Value = "2025-08-07"
$gazDateCell.Style = "Normal"
$wsGaz.Range("B54").Value = 32.175

# ---------------------------------------------------------------------------
# Sheet "CO2": append a new daily row (54) for 2025-08-07.
# ---------------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")
$co2DateCell = $wsCo2.Range("A54")
$co2DateCell.NumberFormat = "@"
$co2DateCell.Value = "2025-08-07"
$co2DateCell.Style = "Normal"
$wsCo2.Range("B54").Value = 71.15000000000001
